# Add new columns I (I0) and J (IF) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy the formatting used by the other header
#     cells (e.g. H1) onto I1 and J1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows (2-40):
#     I column is always 1.
#     J column duplicates the value already present in column H.
$ws.Range("I2:I40").Value = 1

$ws.Range("H2:H40").Copy()
$ws.Range("J2:J40").PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = 0
